{"js": "// Update the title date line and the five filled rows of the 3-digit \u00f7 1-digit\n// division practice table. Each populated table row is replaced in place,\n// cell-by-cell, so per-run formatting (font/size) is preserved.\n\nconst body = context.document.body;\n\n// --- 1) Title paragraph: \"2024-08-16 Friday\" -> \"2024-08-17 Saturday\" ---\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\nparagraphs.items[0].insertText(\"2024-08-17 Saturday\", \"Replace\");\n\n// --- 2) Table cells: replace text of every populated cell, preserving position ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\n\n// Row index -> new values for that row's 5 columns (only the rows that\n// actually contain text are listed; the blank spacer rows are untouched).\nconst rowUpdates = {\n  0: [\"177\u00f76=29, 3\", \"753\u00f79=83, 6\", \"239\u00f78=29, 7\", \"616\u00f77=88, 0\", \"806\u00f79=89, 5\"],\n  4: [\"925\u00f76=154, 1\", \"810\u00f77=115, 5\", \"168\u00f72=84, 0\", \"185\u00f78=23, 1\", \"448\u00f77=64, 0\"],\n  8: [\"433\u00f74=108, 1\", \"748\u00f79=83, 1\", \"313\u00f74=78, 1\", \"815\u00f73=271, 2\", \"801\u00f79=89, 0\"],\n  12: [\"194\u00f76=32, 2\", \"918\u00f75=183, 3\", \"759\u00f76=126, 3\", \"483\u00f79=53, 6\", \"808\u00f78=101, 0\"],\n  16: [\"691\u00f73=230, 1\", \"946\u00f73=315, 1\", \"140\u00f78=17, 4\", \"194\u00f74=48, 2\", \"726\u00f76=121, 0\"],\n};\n\nfor (const rowIndex of Object.keys(rowUpdates)) {\n  const r = Number(rowIndex);\n  const values = rowUpdates[rowIndex];\n  for (let c = 0; c < values.length; c++) {\n    table.getCell(r, c).value = values[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the title date line and the five filled rows of the 3-digit \u00f7 1-digit\n# division practice table. Each populated cell's Range.Text is rewritten in\n# place so per-run formatting (font/size) carried on the existing run stays\n# intact.\n\n$d = $word.ActiveDocument\n\n# --- 1) Title paragraph: \"2024-08-16 Friday\" -> \"2024-08-17 Saturday\" ---\n$d.Paragraphs.Item(1).Range.Text = \"2024-08-17 Saturday\"\n\n# --- 2) Table cells: replace text of every populated cell, preserving position ---\n$t = $d.Tables.Item(1)\n\n# Row -> new values for that row's 5 columns (1-based row/column indices,\n# matching Word's Table.Cell(row, col); only the rows that actually contain\n# text are listed, the blank spacer rows are left untouched).\n$rowUpdates = @{\n    1  = @(\"177\u00f76=29, 3\", \"753\u00f79=83, 6\", \"239\u00f78=29, 7\", \"616\u00f77=88, 0\", \"806\u00f79=89, 5\")\n    5  = @(\"925\u00f76=154, 1\", \"810\u00f77=115, 5\", \"168\u00f72=84, 0\", \"185\u00f78=23, 1\", \"448\u00f77=64, 0\")\n    9  = @(\"433\u00f74=108, 1\", \"748\u00f79=83, 1\", \"313\u00f74=78, 1\", \"815\u00f73=271, 2\", \"801\u00f79=89, 0\")\n    13 = @(\"194\u00f76=32, 2\", \"918\u00f75=183, 3\", \"759\u00f76=126, 3\", \"483\u00f79=53, 6\", \"808\u00f78=101, 0\")\n    17 = @(\"691\u00f73=230, 1\", \"946\u00f73=315, 1\", \"140\u00f78=17, 4\", \"194\u00f74=48, 2\", \"726\u00f76=121, 0\")\n}\n\nforeach ($row in $rowUpdates.Keys) {\n    $values = $rowUpdates[$row]\n    for ($col = 1; $col -le $values.Length; $col++) {\n        $t.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
